# Add the new "cung đối Phúc Đức" combinations following the existing
# two-column (A/B mirrored) layout used throughout this sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target rows for the new entries. Row 4306 and row 4398 are left blank,
# matching the separator-row pattern already present elsewhere in the sheet
# (e.g. the existing gap between rows 4225 and 4227).
$rows = @(
4307, 4308, 4309, 4310, 4311, 4312, 4313, 4314, 4315, 4316, 4317, 4318, 4319, 4320, 4321, 4322, 4323, 4324, 4325, 4326, 4327, 4328, 4329, 4330, 4331, 4332, 4333, 4334, 4335, 4336, 4337, 4338, 4339, 4340, 4341, 4342, 4343, 4344, 4345, 4346, 4347, 4348, 4349, 4350, 4351, 4352, 4353, 4354, 4355, 4356, 4357, 4358, 4359, 4360, 4361, 4362, 4363, 4364, 4365, 4366, 4367, 4368, 4369, 4370, 4371, 4372, 4373, 4374, 4375, 4376, 4377, 4378, 4379, 4380, 4381, 4382, 4383, 4384, 4385, 4386, 4387, 4388, 4389, 4390, 4391, 4392, 4393, 4394, 4395, 4396, 4397, 4399, 4400, 4401, 4402, 4403, 4404, 4405, 4406, 4407, 4408, 4409, 4410, 4411, 4412
)

$values = @(
    "Tử Vi đồng cung Thiên Cơ tại cung đối Phúc Đức",
    "Tử Vi đồng cung Thái Dương tại cung đối Phúc Đức",
    "Tử Vi đồng cung Vũ Khúc tại cung đối Phúc Đức",
    "Tử Vi đồng cung Thiên Đồng tại cung đối Phúc Đức",
    "Tử Vi đồng cung Liêm Trinh tại cung đối Phúc Đức",
    "Tử Vi đồng cung Thiên Phủ tại cung đối Phúc Đức",
    "Tử Vi đồng cung Thái Âm tại cung đối Phúc Đức",
    "Tử Vi đồng cung Tham Lang tại cung đối Phúc Đức",
    "Tử Vi đồng cung Cự Môn tại cung đối Phúc Đức",
    "Tử Vi đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Tử Vi đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Tử Vi đồng cung Thất Sát tại cung đối Phúc Đức",
    "Tử Vi đồng cung Phá Quân tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Thái Dương tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Vũ Khúc tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Thiên Đồng tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Liêm Trinh tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Thiên Phủ tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Thái Âm tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Tham Lang tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Cự Môn tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Thất Sát tại cung đối Phúc Đức",
    "Thiên Cơ đồng cung Phá Quân tại cung đối Phúc Đức",
    "Thái Dương đồng cung Vũ Khúc tại cung đối Phúc Đức",
    "Thái Dương đồng cung Thiên Đồng tại cung đối Phúc Đức",
    "Thái Dương đồng cung Liêm Trinh tại cung đối Phúc Đức",
    "Thái Dương đồng cung Thiên Phủ tại cung đối Phúc Đức",
    "Thái Dương đồng cung Thái Âm tại cung đối Phúc Đức",
    "Thái Dương đồng cung Tham Lang tại cung đối Phúc Đức",
    "Thái Dương đồng cung Cự Môn tại cung đối Phúc Đức",
    "Thái Dương đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Thái Dương đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Thái Dương đồng cung Thất Sát tại cung đối Phúc Đức",
    "Thái Dương đồng cung Phá Quân tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Thiên Đồng tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Liêm Trinh tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Thiên Phủ tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Thái Âm tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Tham Lang tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Cự Môn tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Thất Sát tại cung đối Phúc Đức",
    "Vũ Khúc đồng cung Phá Quân tại cung đối Phúc Đức",
    "Thiên Đồng đồng cung Liêm Trinh tại cung đối Phúc Đức",
    "Thiên Đồng đồng cung Thiên Phủ tại cung đối Phúc Đức",
    "Thiên Đồng đồng cung Thái Âm tại cung đối Phúc Đức",
    "Thiên Đồng đồng cung Tham Lang tại cung đối Phúc Đức",
    "Thiên Đồng đồng cung Cự Môn tại cung đối Phúc Đức",
    "Thiên Đồng đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Thiên Đồng đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Thiên Đồng đồng cung Thất Sát tại cung đối Phúc Đức",
    "Thiên Đồng đồng cung Phá Quân tại cung đối Phúc Đức",
    "Liêm Trinh đồng cung Thiên Phủ tại cung đối Phúc Đức",
    "Liêm Trinh đồng cung Thái Âm tại cung đối Phúc Đức",
    "Liêm Trinh đồng cung Tham Lang tại cung đối Phúc Đức",
    "Liêm Trinh đồng cung Cự Môn tại cung đối Phúc Đức",
    "Liêm Trinh đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Liêm Trinh đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Liêm Trinh đồng cung Thất Sát tại cung đối Phúc Đức",
    "Liêm Trinh đồng cung Phá Quân tại cung đối Phúc Đức",
    "Thiên Phủ đồng cung Thái Âm tại cung đối Phúc Đức",
    "Thiên Phủ đồng cung Tham Lang tại cung đối Phúc Đức",
    "Thiên Phủ đồng cung Cự Môn tại cung đối Phúc Đức",
    "Thiên Phủ đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Thiên Phủ đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Thiên Phủ đồng cung Thất Sát tại cung đối Phúc Đức",
    "Thiên Phủ đồng cung Phá Quân tại cung đối Phúc Đức",
    "Thái Âm đồng cung Tham Lang tại cung đối Phúc Đức",
    "Thái Âm đồng cung Cự Môn tại cung đối Phúc Đức",
    "Thái Âm đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Thái Âm đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Thái Âm đồng cung Thất Sát tại cung đối Phúc Đức",
    "Thái Âm đồng cung Phá Quân tại cung đối Phúc Đức",
    "Tham Lang đồng cung Cự Môn tại cung đối Phúc Đức",
    "Tham Lang đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Tham Lang đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Tham Lang đồng cung Thất Sát tại cung đối Phúc Đức",
    "Tham Lang đồng cung Phá Quân tại cung đối Phúc Đức",
    "Cự Môn đồng cung Thiên Tướng tại cung đối Phúc Đức",
    "Cự Môn đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Cự Môn đồng cung Thất Sát tại cung đối Phúc Đức",
    "Cự Môn đồng cung Phá Quân tại cung đối Phúc Đức",
    "Thiên Tướng đồng cung Thiên Lương tại cung đối Phúc Đức",
    "Thiên Tướng đồng cung Thất Sát tại cung đối Phúc Đức",
    "Thiên Tướng đồng cung Phá Quân tại cung đối Phúc Đức",
    "Thiên Lương đồng cung Thất Sát tại cung đối Phúc Đức",
    "Thiên Lương đồng cung Phá Quân tại cung đối Phúc Đức",
    "Thất Sát đồng cung Phá Quân tại cung đối Phúc Đức",
    "Tử Vi tọa thủ tại cung đối Phúc Đức",
    "Thiên Cơ tọa thủ tại cung đối Phúc Đức",
    "Thái Dương tọa thủ tại cung đối Phúc Đức",
    "Vũ Khúc tọa thủ tại cung đối Phúc Đức",
    "Thiên Đồng tọa thủ tại cung đối Phúc Đức",
    "Liêm Trinh tọa thủ tại cung đối Phúc Đức",
    "Thiên Phủ tọa thủ tại cung đối Phúc Đức",
    "Thái Âm tọa thủ tại cung đối Phúc Đức",
    "Tham Lang tọa thủ tại cung đối Phúc Đức",
    "Cự Môn tọa thủ tại cung đối Phúc Đức",
    "Thiên Tướng tọa thủ tại cung đối Phúc Đức",
    "Thiên Lương tọa thủ tại cung đối Phúc Đức",
    "Thất Sát tọa thủ tại cung đối Phúc Đức",
    "Phá Quân tọa thủ tại cung đối Phúc Đức"
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $v = $values[$i]
    $ws.Cells.Item($r, 1).Value = $v
    $ws.Cells.Item($r, 2).Value = $v
}

# Scroll the view down to the newly added rows and select the last block
# of column B entries, matching where the author left the selection.
try {
    $excel.ActiveWindow.ScrollRow = 4378
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("B4399:B4412").Select()
